$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.745.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = "'2.908.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'528.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.66%  '
$ws.Range("D6").Value = "'145.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.26%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").Value = "'2.917.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("E10").Value = '  -2.64%  '
$ws.Range("D11").Value = "'6.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = "'0.366"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = "'3.416.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = "'60.700.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("E16").Value = '  -2.96%  '
$ws.Range("D17").Value = "'2.914.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = "'11.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").Value = "'363.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.41%  '
$ws.Range("D22").Value = "'6.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = "'64.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = "'0.456"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("D27").Value = "'0.181"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.00%  '
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").Value = "'7.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.32%  '
$ws.Range("E30").Value = '  -5.56%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = "'19.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").Value = "'150.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").Value = "'5.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.22%  '
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("E38").Value = '  -4.41%  '
$ws.Range("D39").Value = "'37.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.80%  '
$ws.Range("E40").Value = '  -2.45%  '
$ws.Range("D41").Value = "'3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").Value = "'2.296.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("D43").Value = "'0.649"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.0583"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").Value = "'20.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.68%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = "'0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'5.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.72%  '
$ws.Range("E48").Value = '  -2.39%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = "'10.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = "'0.0927"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("D51").Value = "'252.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.25%  '
